$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "datos actualizados" timestamp (A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Septiembre de 2020 a las 12:47"

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 6875103
$ws.Range("C4").Value = 507
$ws.Range("D4").Value = 4155655
$ws.Range("E4").Value = 2517229
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 202219

# Row 15: Iran
$ws.Range("A15").Value = "Iran"
$ws.Range("B15").Value = 416198
$ws.Range("C15").Value = 3049
$ws.Range("D15").Value = 355505
$ws.Range("E15").Value = 36741
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 144
$ws.Range("H15").Value = 23952

# Row 16: Francia
$ws.Range("A16").Value = "Francia"
$ws.Range("B16").Value = 415481
$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 90840
$ws.Range("E16").Value = 293546
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 31095

# Row 18: Banglades
$ws.Range("A18").Value = "Banglades"
$ws.Range("B18").Value = 345805
$ws.Range("C18").Value = 1541
$ws.Range("D18").Value = 252335
$ws.Range("E18").Value = 88589
$ws.Range("F18").Value = 0
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = 4881

# Row 49: Bielorrusia
$ws.Range("A49").Value = "Bielorrusia"
$ws.Range("B49").Value = 75230
$ws.Range("C49").Value = 243
$ws.Range("D49").Value = 73098
$ws.Range("E49").Value = 1359
$ws.Range("F49").Value = 0
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 773

# Row 60: Suiza
$ws.Range("A60").Value = "Suiza"
$ws.Range("B60").Value = 49283
$ws.Range("C60").Value = 488
$ws.Range("D60").Value = 39900
$ws.Range("E60").Value = 7340
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 1
$ws.Range("H60").Value = 2043

# Row 61: Argelia
$ws.Range("A61").Value = "Argelia"
$ws.Range("B61").Value = 49194
$ws.Range("C61").Value = 0
$ws.Range("D61").Value = 34675
$ws.Range("E61").Value = 12865
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 0
$ws.Range("H61").Value = 1654

# Row 75: El Salvador
$ws.Range("A75").Value = "El Salvador"
$ws.Range("B75").Value = 27346
$ws.Range("C75").Value = 97
$ws.Range("D75").Value = 20825
$ws.Range("E75").Value = 5717
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 3
$ws.Range("H75").Value = 804

# Row 97: Malasia
$ws.Range("A97").Value = "Malasia"
$ws.Range("B97").Value = 10147
$ws.Range("C97").Value = 95
$ws.Range("D97").Value = 9264
$ws.Range("E97").Value = 754
$ws.Range("F97").Value = 0
$ws.Range("G97").Value = 1
$ws.Range("H97").Value = 129

# Row 98: Namibia
$ws.Range("A98").Value = "Namibia"
$ws.Range("B98").Value = 10078
$ws.Range("C98").Value = 0
$ws.Range("D98").Value = 7685
$ws.Range("E98").Value = 2285
$ws.Range("F98").Value = 0
$ws.Range("G98").Value = 0
$ws.Range("H98").Value = 108

# Row 113: Uganda
$ws.Range("A113").Value = "Uganda"
$ws.Range("B113").Value = 5594
$ws.Range("C113").Value = 214
$ws.Range("D113").Value = 2544
$ws.Range("E113").Value = 2989
$ws.Range("F113").Value = 0
$ws.Range("G113").Value = 1
$ws.Range("H113").Value = 61

# Row 114: Republica de Yibuti
$ws.Range("A114").Value = "Republica de Yibuti"
$ws.Range("B114").Value = 5399
$ws.Range("C114").Value = 0
$ws.Range("D114").Value = 5333
$ws.Range("E114").Value = 5
$ws.Range("F114").Value = 0
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 61

# Row 126: Birmania
$ws.Range("A126").Value = "Birmania"
$ws.Range("B126").Value = 4299
$ws.Range("C126").Value = 256
$ws.Range("D126").Value = 944
$ws.Range("E126").Value = 3287
$ws.Range("F126").Value = 0
$ws.Range("G126").Value = 8
$ws.Range("H126").Value = 68

# Row 145: Malta
$ws.Range("A145").Value = "Malta"
$ws.Range("B145").Value = 2634
$ws.Range("C145").Value = 39
$ws.Range("D145").Value = 1996
$ws.Range("E145").Value = 621
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 1
$ws.Range("H145").Value = 17

# Row 146: Sudan del Sur
$ws.Range("A146").Value = "Sudan del Sur"
$ws.Range("B146").Value = 2599
$ws.Range("C146").Value = 0
$ws.Range("D146").Value = 1290
$ws.Range("E146").Value = 1260
$ws.Range("F146").Value = 0
$ws.Range("G146").Value = 0
$ws.Range("H146").Value = 49

# Row 214: Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0

# Row 215: Montserrat
$ws.Range("A215").Value = "Montserrat"
$ws.Range("B215").Value = 13
$ws.Range("C215").Value = 0
$ws.Range("D215").Value = 12
$ws.Range("E215").Value = 0
$ws.Range("F215").Value = 0
$ws.Range("G215").Value = 0
$ws.Range("H215").Value = 1
